# Replace the synthetic GridLAB-D load names ("load_#_p2ulv#####") in the
# "Load_name" column (F) with the new OCHRE/transactive-load house naming
# scheme ("tl_house_<row-1>"), to reflect switching the co-simulation setup
# to dockerized OCHRE integrated with GridLAB-D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 1; $i -le 40; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 6).Value = "tl_house_$i"
}

# Leave the cursor where the author's last save left it.
$ws.Range("I7").Select() | Out-Null
